# Conserto do erro com o rótulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

function Set-TextLabelKeepingStyle {
    param($ws, $Address, $Text, $StyleSourceAddress)
    $target = $ws.Range($Address)
    # Forcing a "Text" number format before assigning the value keeps Excel
    # from re-interpreting a numeric-looking label (e.g. "2050") as a number.
    $target.NumberFormat = "@"
    $target.Value = $Text
    # The NumberFormat change above forks off a brand new cell style, so
    # copy the original (bold/border/centered) formatting back from a
    # neighboring header cell that still carries the untouched style.
    $ws.Range($StyleSourceAddress).Copy() | Out-Null
    $target.PasteSpecial($xlPasteFormats) | Out-Null
}

# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)")
# and sheet 5 ("Emissoes Totais (MtCO2eq)") all have a mislabeled header in
# E1 that shows the stray number 676.2706852435415 instead of the "2050"
# (or "2041-2050") period label, and sheets 1-4 also have a spurious
# "Total" row at the bottom of their tables that needs to be removed.

# Map of sheet name -> correct E1 label. "Potencia Incremental - SIN(MW)"
# labels its columns as year ranges (2015-2030, 2031-2040, ...), so its
# corrected E1 label follows suit with "2041-2050" instead of "2050".
$tableSheetLabels = [ordered]@{
    "Potencia Acumulada - SIN (MW)"  = "2050"
    "Geracao Periodo Medio (MWMed)"  = "2050"
    "Atendimento a Ponta(MW)"        = "2050"
    "Potencia Incremental - SIN(MW)" = "2041-2050"
}

foreach ($name in $tableSheetLabels.Keys) {
    $ws = $wb.Worksheets.Item($name)
    Set-TextLabelKeepingStyle $ws "E1" $tableSheetLabels[$name] "D1"
    # Remove the bottom "Total" row (row 13).
    $ws.Rows.Item(13).Delete()
}

# "Emissoes Totais (MtCO2eq)" only needs the E1 label fixed; it has no
# Total row.
$wsEmissoes = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabelKeepingStyle $wsEmissoes "E1" "2050" "D1"

# "Custo Total (bilhões de R$)" only has the spurious "Total" row (row 4)
# to remove; its header row has no year-range labels.
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Rows.Item(4).Delete()
